$d = $word.ActiveDocument

function Replace-InScope($range, $oldText, $newText) {
    $ok = $range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $ok) {
        throw "Replace failed for: $oldText"
    }
}

# ------------------------------------------------------------------
# 1) "During sprint 5" -> "During sprint 6" and the rest of that
#    summary paragraph (paragraph 20).
# ------------------------------------------------------------------
$p20 = $d.Paragraphs.Item(20).Range
Replace-InScope $p20 "During sprint 5" "During sprint 6"

$p20 = $d.Paragraphs.Item(20).Range
Replace-InScope $p20 "The hours were entirely project based and we spent 0 time working on documentation" "The hours were mostly spent on the project with minimal time being spent on documentation"

$p20 = $d.Paragraphs.Item(20).Range
Replace-InScope $p20 "The total time we spent 46.5 hours" "The total time we spent 41 hours"

$p20 = $d.Paragraphs.Item(20).Range
Replace-InScope $p20 "finished everything we had to planned with this sprint" "finished everything we had planned for this sprint"

# ------------------------------------------------------------------
# 2) Sprint accomplishments list — headers & descriptions.
# ------------------------------------------------------------------

# Item 1 header / description
$p22 = $d.Paragraphs.Item(22).Range
Replace-InScope $p22 "Checkout Multiple Repos" "Front End And Back End Unit Tests"

$p23 = $d.Paragraphs.Item(23).Range
Replace-InScope $p23 "This was a very important portion of our project. Moving from our prototype where you were only able to check out 1 repo, this allowed you to check out many without resetting the application. Dylan primarily worked on this effort and spent a total of 16 hours." "Nate was in charge of this portion of the project He spent all of his time getting the testing framework stood up and working with our Jenkins build job. He successfully implemented the framework and all tests in 20 hours"

# Item 2 header
$p24 = $d.Paragraphs.Item(24).Range
Replace-InScope $p24 "Research Scoring Algorithm" "Refactor Scoring Algorithm"

# Item 2 description — needs a gramStart/gramEnd proofErr pair around "was",
# so rebuild the whole paragraph's XML precisely.
$p25 = $d.Paragraphs.Item(25).Range
$xml25 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Spencer and Dylan were primary contributors to this task. They spent a combined </w:t></w:r><w:r><w:t>18 hours on this task</w:t></w:r><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>was</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> successful in implementing the time based metric of the algorithm.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p25.InsertXML($xml25)

# Item 3 header
$p26 = $d.Paragraphs.Item(26).Range
Replace-InScope $p26 "Store Repo Scores and Meta Data" "Update Documentation"

# Item 3 description
$p27 = $d.Paragraphs.Item(27).Range
Replace-InScope $p27 "This effort was primary spear headed by Spencer. This allowed our algorithm to only run once and then recall the information from the database when the information is called for. Spencer spent a total of 11 hours on this." "Dylan was in charge of completing this effort. He spent 3 hours updating the Project Plan and Quality Plan to reflect our current project."

# Item 4 header
$p28 = $d.Paragraphs.Item(28).Range
Replace-InScope $p28 "Research and POC of Testing of Framework" "Research Scoring Algorithm"

# Item 4 description — needs a spellStart/spellEnd proofErr pair around "Crk",
# so rebuild the whole paragraph's XML precisely.
$p29 = $d.Paragraphs.Item(29).Range
$xml29 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">This was an effort that is ongoing as we progress with our project. Everyone contributes equally and shares their thoughts on how to better the project. We also consult closely to Dr. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Crk</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> on what we need to change about our algorithm and methods on how to test it.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p29.InsertXML($xml29)

# ------------------------------------------------------------------
# 3) Remove the blank ListParagraph spacer paragraph that used to sit
#    right after the last description (now paragraph 30).
# ------------------------------------------------------------------
$d.Paragraphs.Item(30).Range.Delete()

# ------------------------------------------------------------------
# 4) Client update paragraph — the "He wants us..." sentence becomes a
#    much longer passage, including a second spellStart/spellEnd "Crk".
#    (this is now paragraph 31 after the delete above)
# ------------------------------------------------------------------
$p31 = $d.Paragraphs.Item(31).Range
$xml31 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">After the </w:t></w:r><w:r><w:t>completion</w:t></w:r><w:r><w:t xml:space="preserve"> of these tasks, we updated our client Dr. Igor Crk. </w:t></w:r><w:r><w:t xml:space="preserve">He was very pleased with our team being able to produce the algorithm in his original specification. </w:t></w:r><w:r><w:t>A</w:t></w:r><w:r><w:t>s of right now</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> we </w:t></w:r><w:r><w:t xml:space="preserve">are </w:t></w:r><w:r><w:t xml:space="preserve">100% finished </w:t></w:r><w:r><w:t xml:space="preserve">with </w:t></w:r><w:r><w:t xml:space="preserve">everything Dr. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Crk</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> want</w:t></w:r><w:r><w:t>ed</w:t></w:r><w:r><w:t xml:space="preserve"> us to implement. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p31.InsertXML($xml31)

# ------------------------------------------------------------------
# 5) Risks / consensus paragraph (now paragraph 32).
# ------------------------------------------------------------------
$p32 = $d.Paragraphs.Item(32).Range
Replace-InScope $p32 "The major risks posed this sprint was getting the testing framework set up and tests wrote. It has been found to be very difficult to orchestrate meaningful tests with how our project is currently set up. Much time and refactoring is needed by all members to master this art. " "The major risks posed for this sprint was getting the testing framework set up and tests wrote. It has been found to be very difficult to orchestrate meaningful tests with how our project is currently set up. Nate spent most the time refactoring the entirety of the front end to mitigate this risk. "

$p32 = $d.Paragraphs.Item(32).Range
Replace-InScope $p32 "still very optimistic. Having a fully working application allows us to really focus on molding it into something meaningful." "still great due to the fact we have finished everything that our client wanted use to finish. Being able to move forward and add things we feel will benefit the application make the project seem like less work and is motivating."

# ------------------------------------------------------------------
# 6) Move the _GoBack bookmark from the Chart 13 paragraph down to the
#    blank paragraph right after the risks/consensus paragraph. Adding
#    a new _GoBack bookmark implicitly relocates the (singleton)
#    existing one.
# ------------------------------------------------------------------
$spacerPara = $d.Paragraphs.Item($d.Paragraphs.Item(32).Index + 1)
$d.Bookmarks.Add("_GoBack", $spacerPara.Range)

# ------------------------------------------------------------------
# 7) Mark the Signature.jpg picture run as NoProofing (adds <w:noProof/>).
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shape = $d.InlineShapes.Item($i)
    $before = $d.Range([Math]::Max(0, $shape.Range.Start - 40), $shape.Range.Start)
    if ($before.Text -like "*Reinhardt*__*") {
        $shape.Range.NoProofing = $true
    }
}

# ------------------------------------------------------------------
# 8) Footer page-number field cache: "2" -> "4".
# ------------------------------------------------------------------
$footer = $d.Sections.Item(1).Footers.Item(1)
$footer.Range.Characters.Item(1).Text = "4"
